# Update "want-to-go" counts (column F) for several rows across the
# 展览 (Exhibition), 本地生活 (Local life) and 全部类型 (All types) sheets,
# matching the regenerated gh-pages data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet order in the workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wsExhibition = $wb.Worksheets.Item(1)
$wsLocalLife  = $wb.Worksheets.Item(3)
$wsAllTypes   = $wb.Worksheets.Item(4)

# --- 展览 (Exhibition) sheet ---
$wsExhibition.Range("F2").Value  = 1568
$wsExhibition.Range("F6").Value  = 8008
$wsExhibition.Range("F10").Value = 5780
$wsExhibition.Range("F14").Value = 8156
$wsExhibition.Range("F15").Value = 9513
$wsExhibition.Range("F17").Value = 951
$wsExhibition.Range("F20").Value = 287
$wsExhibition.Range("F24").Value = 1226
$wsExhibition.Range("F26").Value = 1734
$wsExhibition.Range("F28").Value = 998
$wsExhibition.Range("F29").Value = 118
$wsExhibition.Range("F32").Value = 496
$wsExhibition.Range("F33").Value = 2396
$wsExhibition.Range("F44").Value = 445
$wsExhibition.Range("F46").Value = 19

# --- 本地生活 (Local life) sheet ---
$wsLocalLife.Range("F2").Value = 5433

# --- 全部类型 (All types) sheet ---
$wsAllTypes.Range("F2").Value  = 1568
$wsAllTypes.Range("F10").Value = 5780
$wsAllTypes.Range("F12").Value = 8156
$wsAllTypes.Range("F13").Value = 9513
$wsAllTypes.Range("F16").Value = 951
$wsAllTypes.Range("F18").Value = 287
$wsAllTypes.Range("F23").Value = 1226
$wsAllTypes.Range("F26").Value = 998
$wsAllTypes.Range("F27").Value = 119
$wsAllTypes.Range("F30").Value = 496
$wsAllTypes.Range("F31").Value = 2396
$wsAllTypes.Range("F43").Value = 445
$wsAllTypes.Range("F45").Value = 19

$wb.Save()
